# Append two new contact-form submissions to the sheet (rows 10 and 11),
# mirroring the existing header/data layout (A:Name, B:Email, C:Contact,
# D:Brand, E:Service, F:Timestamp).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - karthika
$ws.Range("A10").Value = "karthika"
$ws.Range("B10").Value = "karthi@gmail.com"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "895"
$ws.Range("D10").Value = "hh"
$ws.Range("E10").Value = "btl"
$ws.Range("F10").Value = "6/28/2025, 4:28:29 PM"

# Row 11 - xnbZ Xnb X
$ws.Range("A11").Value = "xnbZ Xnb X"
$ws.Range("B11").Value = "madhumitha.24mca@kct.ac.in"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "09865856968"
$ws.Range("D11").Value = "ss"
$ws.Range("E11").Value = "branding"
$ws.Range("F11").Value = "6/28/2025, 5:05:03 PM"
